$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("!!_Table of contents")
$ws1.Unprotect()
$ws1.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.9' date='2020-04-27 01:04:59'"
$ws1.Range("A2").Value = "!!ObjTables type='TableOfContents' tableFormat='row' description='Table of contents' date='2020-04-27 01:04:59' objTablesVersion='0.0.9'"
$ws1.Protect($null, $true, $true, $true, $true)

$ws2 = $wb.Worksheets.Item("!!_Schema")
$ws2.Unprotect()
$ws2.Range("A1").Value = "!!ObjTables type='Schema' tableFormat='row' description='Table/model and column/attribute definitions' date='2020-04-27 01:04:59' objTablesVersion='0.0.9'"
$ws2.Protect($null, $true, $true, $true, $true)

$ws3 = $wb.Worksheets.Item("!!Company")
$ws3.Unprotect()
$ws3.Range("A1").Value = "!!ObjTables type='Data' tableFormat='column' class='Company' name='Companies' date='2020-04-27 01:04:59' objTablesVersion='0.0.9'"
$ws3.Protect($null, $true, $true, $true, $true)

$ws4 = $wb.Worksheets.Item("!!People")
$ws4.Unprotect()
$ws4.Range("A1").Value = "!!ObjTables type='Data' tableFormat='row' class='Person' name='People' date='2020-04-27 01:04:59' objTablesVersion='0.0.9'"
$ws4.Protect($null, $true, $true, $true, $true)
